$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# The "Periodo Mora" column (E16:E26) is re-populated with a new/reversed
# set of periods, and the "Valor Mora" amounts in F16 and F26 are swapped
# as part of the database update.

$ws.Range("E16").Value = "1801"
$ws.Range("E17").Value = "1712"
$ws.Range("E18").Value = "1711"
$ws.Range("E19").Value = "1710"
$ws.Range("E20").Value = "1709"
$ws.Range("E21").Value = "1708"
$ws.Range("E22").Value = "1707"
$ws.Range("E23").Value = "1706"
$ws.Range("E24").Value = "1705"
$ws.Range("E25").Value = "1704"
$ws.Range("E26").Value = "1703"

$ws.Range("F16").Value = 28526
$ws.Range("F26").Value = 29509
